$d = $word.ActiveDocument

# Locate the unique anchor paragraph text ("...perform this task.  Use the
# -data option ... ascii text.") so we operate on the correct occurrence
# of "ascii text." (the phrase also doubles as a unique search key on its
# own, but we verify context first to be safe).
$anchor = $d.Content.Duplicate
$anchorFound = $anchor.Find.Execute("perform this task", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $anchorFound) {
    throw "Could not locate the 'perform this task' anchor paragraph"
}

$tail = $d.Range($anchor.End, $d.Content.End)
$tailFound = $tail.Find.Execute("ascii text.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $tailFound) {
    throw "Could not locate the 'ascii text.' phrase to extend"
}

# $tail now spans exactly the phrase "ascii text." at the end of the
# paragraph run. Replace it in place (preserves the run's character
# formatting) with the merged sentence plus the new explanatory text
# describing the nping flags needed for the TCP session hijacking task.
$tail.Text = "ascii text. You will also want to provide the psh and ack flags, and ack the previous packet in your spoofed packet. Your goal is to use a spoofed packet to hijack a telnet session and delete the file on the server at ~/documents/delete-this.txt.  Note that if you use your telnet session to delete that file, e.g., to observe the protocol in wireshark, then you must recreate that file so it can be deleted in a hijacked session."
